$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column A (ID column) for rows 2 and 4
$ws.Range("A2").Value = 40412
$ws.Range("A4").Value = 40416

# Move the active selection to F7, matching the final cursor position left by the author
$ws.Range("F7").Select()
